$d = $word.ActiveDocument

# 1. Add a collapsed "_GoBack" bookmark right after the "撤销操作" heading text.
$rng = $d.Content
$rng.Find.Execute("撤销操作")
$d.Bookmarks.Add("_GoBack", $rng)

# 2. "提交完了" -> "提交"
$d.Content.Find.Execute("提交完了", $true, $false, $false, $false, $false, $true, 1, $false, "提交", 2)

# 3. "后，想" -> "(commit)后，想"
$d.Content.Find.Execute("后，想", $true, $false, $false, $false, $false, $true, 1, $false, "(commit)后，想", 2)

# 4. Merge "跳过使用暂" + (old _GoBack bookmark) + "存区" into a single run
#    "跳过使用暂存区", removing the old bookmark in the process, while leaving
#    the following "d" and "域" runs untouched.
$d.Content.Find.Execute("跳过使用暂存区", $true, $false, $false, $false, $false, $true, 1, $false, "跳过使用暂存区", 2)

# 5. "asd" -> "asddd"
$d.Content.Find.Execute("asd", $true, $false, $false, $false, $false, $true, 1, $false, "asddd", 2)
